$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 2.44
$ws.Range("G2").Value = 3.1
$ws.Range("H2").Value = 2.5
$ws.Range("I2").Value = 3.6
$ws.Range("K2").Value = 5.4
$ws.Range("P2").Value = 1.71
# Row 3
$ws.Range("G3").Value = 3.45
$ws.Range("I3").Value = 3.2
$ws.Range("K3").Value = 5.1
$ws.Range("N3").Value = 1.88
$ws.Range("P3").Value = 1.88
$ws.Range("Q3").Value = 1.67
$ws.Range("R3").Value = 1.35
$ws.Range("S3").Value = 2.62
# Row 4
$ws.Range("S4").Value = 1.84
# Row 5
$ws.Range("F5").Value = 3.5
$ws.Range("H5").Value = 2.1
$ws.Range("I5").Value = 2.12
$ws.Range("K5").Value = 4.1
$ws.Range("X5").Value = 40
# Row 6
$ws.Range("G6").Value = 3.35
$ws.Range("H6").Value = 2.2
$ws.Range("K6").Value = 4.1
$ws.Range("Q6").Value = 1.53
$ws.Range("R6").Value = 1.71
$ws.Range("T6").Value = 1.51
$ws.Range("X6").Value = 28
$ws.Range("AA6").Value = 48
$ws.Range("AB6").Value = 21
$ws.Range("AD6").Value = 12
$ws.Range("AN6").Value = 19
$ws.Range("AO6").Value = 9.800000000000001
# Row 12
$ws.Range("F12").Value = 1.85
$ws.Range("G12").Value = 1.86
$ws.Range("I12").Value = 5.1
$ws.Range("J12").Value = 3.8
$ws.Range("K12").Value = 3.9
$ws.Range("N12").Value = 3.7
$ws.Range("O12").Value = 1.33
$ws.Range("P12").Value = 1.92
$ws.Range("AN12").Value = 12.5
# Row 13
$ws.Range("G13").Value = 9.800000000000001
$ws.Range("K13").Value = 6.8
$ws.Range("S13").Value = 1.95
$ws.Range("Z13").Value = 12.5
# Row 14
$ws.Range("F14").Value = 3.9
$ws.Range("G14").Value = 4
$ws.Range("I14").Value = 2.02
# Row 15
$ws.Range("F15").Value = 1.89
$ws.Range("G15").Value = 1.95
$ws.Range("J15").Value = 3.65
$ws.Range("K15").Value = 3.8
$ws.Range("Q15").Value = 2.02
# Row 16
$ws.Range("K16").Value = 5.6
$ws.Range("AB16").Value = 15.5
$ws.Range("AJ16").Value = 15.5
# Row 17
$ws.Range("I17").Value = 23
$ws.Range("J17").Value = 9.199999999999999
$ws.Range("K17").Value = 9.6
$ws.Range("R17").Value = 1.96
$ws.Range("S17").Value = 1.93
# Row 18
$ws.Range("H18").Value = 22
$ws.Range("I18").Value = 24
$ws.Range("K18").Value = 12.5
$ws.Range("R18").Value = 2.4
